# Add word-boundary wrapping to the disease search-term regexes on the
# "List" sheet (column D, rows 2-11): turn each raw pattern into
#   (?:[^a-zA-Z]|\b)(<original pattern>)(?:[^a-zA-Z]|\b)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List")
$ws.Activate()

$lastRow = $ws.Cells.Item($ws.Rows.Count(), 1).End(-4162).Row()

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $pattern = $cell.Value()
    if ($pattern) {
        $cell.Value = "(?:[^a-zA-Z]|\b)(" + $pattern + ")(?:[^a-zA-Z]|\b)"
    }
}

# Restore the selection left by the editing session.
$ws.Range("D5").Select()
